$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value2 = 647.63635
$ws.Range("I12").Value2 = 687.4
$ws.Range("J12").Value2 = 250
$ws.Range("K12").Value2 = 687.4
$ws.Range("L12").Value2 = 250
$ws.Range("M12").Value2 = -517.4
$ws.Range("N12").Value2 = -590
$ws.Range("H32").Value2 = 12503663
$ws.Range("I32").Value2 = 1200
$ws.Range("J32").Value2 = 14289729
$ws.Range("K32").Value2 = 1200
$ws.Range("L32").Value2 = 14289729
$ws.Range("M32").Value2 = -874
$ws.Range("N32").Value2 = -14290381
$ws.Range("H43").Value2 = 4953.4614
$ws.Range("I43").Value2 = 4779.4
$ws.Range("J43").Value2 = 5062.25
$ws.Range("K43").Value2 = 4779.4
$ws.Range("L43").Value2 = 5062.25
$ws.Range("M43").Value2 = -4710.4
$ws.Range("N43").Value2 = -5200.25
$ws.Range("H53").Value2 = 297.36365
$ws.Range("I53").Value2 = 115.5
$ws.Range("J53").Value2 = 365.5625
$ws.Range("K53").Value2 = 115.5
$ws.Range("L53").Value2 = 365.5625
$ws.Range("M53").Value2 = 521.5
$ws.Range("N53").Value2 = -1639.5625
$ws.Range("H70").Value2 = 1760.2
$ws.Range("I70").Value2 = 1100.5
$ws.Range("J70").Value2 = 2200
$ws.Range("K70").Value2 = 3301.5
$ws.Range("L70").Value2 = 6600
$ws.Range("M70").Value2 = -3031.5
$ws.Range("N70").Value2 = -7140
$ws.Range("H73").Value2 = 1760.2
$ws.Range("I73").Value2 = 1100.5
$ws.Range("J73").Value2 = 2200
$ws.Range("K73").Value2 = 3301.5
$ws.Range("L73").Value2 = 6600
$ws.Range("M73").Value2 = -2365.5
$ws.Range("N73").Value2 = -8472
$ws.Range("H86").Value2 = 3610.05
$ws.Range("I86").Value2 = 2992.375
$ws.Range("J86").Value2 = 4021.8333
$ws.Range("K86").Value2 = 2992.375
$ws.Range("L86").Value2 = 4021.8333
$ws.Range("M86").Value2 = -1869.375
$ws.Range("N86").Value2 = -6267.8333
$ws.Range("H89").Value2 = 3610.05
$ws.Range("I89").Value2 = 2992.375
$ws.Range("J89").Value2 = 4021.8333
$ws.Range("K89").Value2 = 14961.875
$ws.Range("L89").Value2 = 20109.1665
$ws.Range("M89").Value2 = -9345.875
$ws.Range("N89").Value2 = -31341.1665
$ws.Range("H98").Value2 = 0
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 0
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 0
$ws.Range("M98").ClearContents()
$ws.Range("H122").Value2 = 0
$ws.Range("I122").Value2 = 0
$ws.Range("J122").Value2 = 0
$ws.Range("K122").Value2 = 0
$ws.Range("L122").Value2 = 0
$ws.Range("M122").ClearContents()
$ws.Range("H137").Value2 = 6514.32
$ws.Range("I137").Value2 = 8759.083000000001
$ws.Range("J137").Value2 = 4442.231
$ws.Range("K137").Value2 = 26277.249
$ws.Range("L137").Value2 = 13326.693
$ws.Range("M137").Value2 = -23727.249
$ws.Range("N137").Value2 = -18426.693

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 3886.8472
$ws.Range("I32").Value2 = 4391.0166
$ws.Range("J32").Value2 = 1366
$ws.Range("K32").Value2 = 4391.0166
$ws.Range("L32").Value2 = 1366
$ws.Range("M32").Value2 = -4104.0166
$ws.Range("N32").Value2 = -1940
$ws.Range("H61").Value2 = 4670.769
$ws.Range("I61").Value2 = 4572.3
$ws.Range("J61").Value2 = 4999
$ws.Range("K61").Value2 = 4572.3
$ws.Range("L61").Value2 = 4999
$ws.Range("M61").Value2 = -4360.3
$ws.Range("N61").Value2 = -5423
$ws.Range("H69").Value2 = 80000
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 80000
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 80000
$ws.Range("N69").Value2 = -81498
$ws.Range("H72").Value2 = 80000
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 80000
$ws.Range("K72").Value2 = 0
$ws.Range("L72").Value2 = 240000
$ws.Range("N72").Value2 = -247488
$ws.Range("H74").Value2 = 1301.6818
$ws.Range("I74").Value2 = 1342.9375
$ws.Range("J74").Value2 = 1191.6666
$ws.Range("K74").Value2 = 1342.9375
$ws.Range("L74").Value2 = 1191.6666
$ws.Range("M74").Value2 = -468.9375
$ws.Range("N74").Value2 = -2939.6666
$ws.Range("H77").Value2 = 1301.6818
$ws.Range("I77").Value2 = 1342.9375
$ws.Range("J77").Value2 = 1191.6666
$ws.Range("K77").Value2 = 6714.6875
$ws.Range("L77").Value2 = 5958.333000000001
$ws.Range("M77").Value2 = -2346.6875
$ws.Range("N77").Value2 = -14694.333
$ws.Range("H97").Value2 = 2573.3125
$ws.Range("I97").Value2 = 1210
$ws.Range("J97").Value2 = 5572.6
$ws.Range("K97").Value2 = 1210
$ws.Range("L97").Value2 = 5572.6
$ws.Range("M97").Value2 = -714
$ws.Range("N97").Value2 = -6564.6
$ws.Range("H111").Value2 = 0
$ws.Range("I111").Value2 = 0
$ws.Range("J111").Value2 = 0
$ws.Range("K111").Value2 = 0
$ws.Range("L111").ClearContents()
$ws.Range("N111").Value2 = 0
$ws.Range("H122").Value2 = 5047.68
$ws.Range("I122").Value2 = 4621.8335
$ws.Range("J122").Value2 = 6142.7144
$ws.Range("K122").Value2 = 13865.5005
$ws.Range("L122").Value2 = 18428.1432
$ws.Range("M122").Value2 = -11415.5005
$ws.Range("N122").Value2 = -23328.1432
$ws.Range("H136").Value2 = 4670.769
$ws.Range("I136").Value2 = 4572.3
$ws.Range("J136").Value2 = 4999
$ws.Range("K136").Value2 = 13716.9
$ws.Range("L136").Value2 = 14997
$ws.Range("M136").Value2 = -11166.9
$ws.Range("N136").Value2 = -20097

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H10").Value2 = 6375
$ws.Range("I10").Value2 = 2750
$ws.Range("J10").Value2 = 10000
$ws.Range("K10").Value2 = 2750
$ws.Range("L10").Value2 = 10000
$ws.Range("M10").Value2 = -2610
$ws.Range("N10").Value2 = -10280
$ws.Range("H94").Value2 = 1883.8077
$ws.Range("I94").Value2 = 734.3158
$ws.Range("J94").Value2 = 5003.857
$ws.Range("K94").Value2 = 734.3158
$ws.Range("L94").Value2 = 5003.857
$ws.Range("M94").Value2 = -283.3158
$ws.Range("N94").Value2 = -5905.857
$ws.Range("H99").Value2 = 47731.637
$ws.Range("I99").Value2 = 1214.8
$ws.Range("J99").Value2 = 86495.664
$ws.Range("K99").Value2 = 1214.8
$ws.Range("L99").Value2 = 86495.664
$ws.Range("M99").Value2 = 283.2
$ws.Range("N99").Value2 = -89491.664
$ws.Range("H107").Value2 = 8723
$ws.Range("I107").Value2 = 12063.375
$ws.Range("J107").Value2 = 6050.7
$ws.Range("K107").Value2 = 12063.375
$ws.Range("L107").Value2 = 6050.7
$ws.Range("M107").Value2 = -10143.375
$ws.Range("N107").Value2 = -9890.700000000001

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value2 = 25
$ws.Range("I4").Value2 = 25
$ws.Range("J4").Value2 = 0
$ws.Range("K4").Value2 = 25
$ws.Range("L4").Value2 = 0
$ws.Range("M4").Value2 = 87
$ws.Range("H31").Value2 = 2248.88
$ws.Range("I31").Value2 = 1774.8
$ws.Range("J31").Value2 = 2960
$ws.Range("K31").Value2 = 1774.8
$ws.Range("L31").Value2 = 2960
$ws.Range("M31").Value2 = -1479.8
$ws.Range("N31").Value2 = -3550
$ws.Range("H34").Value2 = 2248.88
$ws.Range("I34").Value2 = 1774.8
$ws.Range("J34").Value2 = 2960
$ws.Range("K34").Value2 = 1774.8
$ws.Range("L34").Value2 = 2960
$ws.Range("M34").Value2 = -1572.8
$ws.Range("N34").Value2 = -3364
$ws.Range("H44").Value2 = 0
$ws.Range("I44").Value2 = 0
$ws.Range("J44").Value2 = 0
$ws.Range("K44").Value2 = 0
$ws.Range("L44").Value2 = 0
$ws.Range("M44").ClearContents()
$ws.Range("H93").Value2 = 38142.715
$ws.Range("I93").Value2 = 4250
$ws.Range("J93").Value2 = 83333
$ws.Range("K93").Value2 = 4250
$ws.Range("L93").Value2 = 83333
$ws.Range("M93").Value2 = -2378
$ws.Range("N93").Value2 = -87077
$ws.Range("H107").Value2 = 4273.2
$ws.Range("I107").Value2 = 719.25
$ws.Range("J107").Value2 = 11381.1
$ws.Range("K107").Value2 = 719.25
$ws.Range("L107").Value2 = 11381.1
$ws.Range("M107").Value2 = 1200.75
$ws.Range("N107").Value2 = -15221.1
$ws.Range("H112").Value2 = 29769.23
$ws.Range("I112").Value2 = 0
$ws.Range("J112").Value2 = 29769.23
$ws.Range("K112").Value2 = 0
$ws.Range("L112").Value2 = 29769.23
$ws.Range("N112").Value2 = -32723.23

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value2 = 166666720
$ws.Range("I44").Value2 = 83
$ws.Range("J44").Value2 = 333333340
$ws.Range("K44").Value2 = 249
$ws.Range("L44").Value2 = 1000000020
$ws.Range("M44").Value2 = 149
$ws.Range("N44").Value2 = -1000000816
$ws.Range("H46").Value2 = 91422.82000000001
$ws.Range("I46").Value2 = 333460.34
$ws.Range("J46").Value2 = 658.75
$ws.Range("K46").Value2 = 1000381.02
$ws.Range("L46").Value2 = 1976.25
$ws.Range("M46").Value2 = -1000290.02
$ws.Range("N46").Value2 = -2158.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value2 = 4485.9585
$ws.Range("I97").Value2 = 341.44446
$ws.Range("J97").Value2 = 16919.5
$ws.Range("K97").Value2 = 341.44446
$ws.Range("L97").Value2 = 16919.5
$ws.Range("M97").Value2 = 154.55554
$ws.Range("N97").Value2 = -17911.5
$ws.Range("H107").Value2 = 579.8889
$ws.Range("I107").Value2 = 417.14285
$ws.Range("J107").Value2 = 1149.5
$ws.Range("K107").Value2 = 417.14285
$ws.Range("L107").Value2 = 1149.5
$ws.Range("M107").Value2 = 1502.85715
$ws.Range("N107").Value2 = -4989.5
$ws.Range("H137").Value2 = 50000
$ws.Range("I137").Value2 = 0
$ws.Range("J137").Value2 = 50000
$ws.Range("K137").Value2 = 0
$ws.Range("L137").Value2 = 50000
$ws.Range("N137").Value2 = -60200

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value2 = 430.7143
$ws.Range("I10").Value2 = 430.7143
$ws.Range("J10").Value2 = 0
$ws.Range("K10").Value2 = 430.7143
$ws.Range("L10").Value2 = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value2 = -290.7143
$ws.Range("H12").Value2 = 0
$ws.Range("I12").Value2 = 0
$ws.Range("J12").Value2 = 0
$ws.Range("K12").Value2 = 0
$ws.Range("L12").ClearContents()
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value2 = 0
$ws.Range("H55").Value2 = 318.08334
$ws.Range("I55").Value2 = 375.375
$ws.Range("J55").Value2 = 203.5
$ws.Range("K55").Value2 = 375.375
$ws.Range("L55").Value2 = 203.5
$ws.Range("M55").Value2 = -202.375
$ws.Range("N55").Value2 = -549.5
$ws.Range("H122").Value2 = 58827980
$ws.Range("I122").Value2 = 142859490
$ws.Range("J122").Value2 = 5919.2
$ws.Range("K122").Value2 = 428578470
$ws.Range("L122").Value2 = 17757.6
$ws.Range("M122").Value2 = -428576020
$ws.Range("N122").Value2 = -22657.6
$ws.Range("H132").Value2 = 2519.8333
$ws.Range("I132").Value2 = 2519.8333
$ws.Range("J132").Value2 = 0
$ws.Range("K132").Value2 = 7559.499899999999
$ws.Range("L132").Value2 = 0
$ws.Range("M132").Value2 = -5029.499899999999
$ws.Range("H136").Value2 = 3415.3125
$ws.Range("I136").Value2 = 2617.7273
$ws.Range("J136").Value2 = 5170
$ws.Range("K136").Value2 = 7853.1819
$ws.Range("L136").Value2 = 15510
$ws.Range("M136").Value2 = -5303.1819
$ws.Range("N136").Value2 = -20610

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value2 = 5489
$ws.Range("I41").Value2 = 0
$ws.Range("J41").Value2 = 5489
$ws.Range("K41").Value2 = 0
$ws.Range("L41").Value2 = 5489
$ws.Range("N41").Value2 = -6269
$ws.Range("H45").Value2 = 25000
$ws.Range("I45").Value2 = 0
$ws.Range("J45").Value2 = 25000
$ws.Range("K45").Value2 = 0
$ws.Range("L45").Value2 = 25000
$ws.Range("N45").Value2 = -25982
$ws.Range("H126").Value2 = 9456.429
$ws.Range("I126").Value2 = 10282.5
$ws.Range("J126").Value2 = 4500
$ws.Range("K126").Value2 = 30847.5
$ws.Range("L126").Value2 = 13500
$ws.Range("M126").Value2 = -28377.5
$ws.Range("N126").Value2 = -18440
$ws.Range("H132").Value2 = 3005.25
$ws.Range("I132").Value2 = 2858.111
$ws.Range("J132").Value2 = 3446.6667
$ws.Range("K132").Value2 = 8574.332999999999
$ws.Range("L132").Value2 = 10340.0001
$ws.Range("M132").Value2 = -6044.332999999999
$ws.Range("N132").Value2 = -15400.0001
$ws.Range("H138").Value2 = 59992.5
$ws.Range("I138").Value2 = 0
$ws.Range("J138").Value2 = 59992.5
$ws.Range("K138").Value2 = 0
$ws.Range("L138").Value2 = 59992.5
$ws.Range("N138").Value2 = -70272.5
